# Zeiterfassung.xlsx update
# - New time-tracking entry on Tabelle1 (row 4): 04.05.2014, Roman, Development,
#   240 minutes, with a status comment about included libs / contact insert & select.
# - Move the active selection to E5 (was E8) to reflect the newly entered row.
# Formulas in G2/G6 (SUMIFS-based totals) recalc automatically off the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New entry: 04.05.2014, Roman, Development, 240 min, status note
$ws.Range("A4").Value = 41763
$ws.Range("B4").Value = "Roman"
$ws.Range("C4").Value = "Development"
$ws.Range("D4").Value = 240
$ws.Range("E4").Value = "Libs inkludiert, Programm lauffähig gemacht, Insert von Kontakten funktioniert, Selecten funktioniert serverseitig"

# Reflect the updated selection/active cell on the sheet
$ws.Range("E5").Select()
